$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.061.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.182.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.180.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("E11").Value = "  -7.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.512"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.703.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.114.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.178.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.729"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.34%  "
$ws.Range("E30").Value = "  +4.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("E34").Value = "  -4.67%  "
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "512.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0883"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.98%  "
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0709"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.95%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.127"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.832.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.117"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.18%  "
